$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '42.312.86'
Set-TextValue 'E2' '  -0.71%  '
Set-TextValue 'D3' '2.273.84'
Set-TextValue 'E3' '  -0.88%  '
Set-TextValue 'E4' '  -0.07%  '
Set-TextValue 'D5' '308.89'
Set-TextValue 'E5' '  +0.19%  '
Set-TextValue 'D6' '97.58'
Set-TextValue 'E6' '  -1.35%  '
Set-TextValue 'D7' '0.527'
Set-TextValue 'E7' '  -1.01%  '
Set-TextValue 'E8' '  +0.03%  '
Set-TextValue 'D9' '0.489'
Set-TextValue 'D10' '35.03'
Set-TextValue 'E10' '  -4.29%  '
Set-TextValue 'D11' '0.0815'
Set-TextValue 'E11' '  +0.89%  '
Set-TextValue 'E12' '  +0.97%  '
Set-TextValue 'D13' '6.83'
Set-TextValue 'E13' '  +0.99%  '
Set-TextValue 'D14' '2.627.88'
Set-TextValue 'E14' '  -0.80%  '
Set-TextValue 'D15' '14.62'
Set-TextValue 'E15' '  -0.17%  '
Set-TextValue 'D16' '2.260.42'
Set-TextValue 'E16' '  -2.39%  '
Set-TextValue 'D17' '0.788'
Set-TextValue 'E17' '  -2.65%  '
Set-TextValue 'D18' '42.191.05'
Set-TextValue 'E18' '  -0.76%  '
Set-TextValue 'D19' '12.27'
Set-TextValue 'E19' '  -4.46%  '
Set-TextValue 'D20' '0.0₃0907'
Set-TextValue 'E20' '  -1.32%  '
Set-TextValue 'D21' '5.97'
Set-TextValue 'E21' '  -1.55%  '
Set-TextValue 'D22' '67.61'
Set-TextValue 'E22' '  -0.72%  '
Set-TextValue 'D23' '236.71'
Set-TextValue 'E23' '  -2.81%  '
Set-TextValue 'E24' '  -1.03%  '
Set-TextValue 'D25' '1.97'
Set-TextValue 'E25' '  +0.04%  '
Set-TextValue 'E26' '  -0.12%  '
Set-TextValue 'D27' '23.57'
Set-TextValue 'E27' '  -1.98%  '
Set-TextValue 'D28' '37.35'
Set-TextValue 'E28' '  -3.01%  '
Set-TextValue 'D29' '9.58'
Set-TextValue 'E29' '  -0.63%  '
Set-TextValue 'E30' '  -0.02%  '
Set-TextValue 'D31' '163.52'
Set-TextValue 'E31' '  +1.67%  '
Set-TextValue 'D32' '5.25'
Set-TextValue 'E32' '  -1.41%  '
Set-TextValue 'E33' '  +0.08%  '
Set-TextValue 'E34' '  -2.00%  '
Set-TextValue 'D35' '17.70'
Set-TextValue 'E35' '  +1.30%  '
Set-TextValue 'D36' '0.0734'
Set-TextValue 'E36' '  -2.44%  '
Set-TextValue 'E37' '  -0.20%  '
Set-TextValue 'E38' '  -3.83%  '
Set-TextValue 'D39' '1.82'
Set-TextValue 'E39' '  -3.41%  '
Set-TextValue 'E40' '  -0.94%  '
Set-TextValue 'D41' '4.17'
Set-TextValue 'E41' '  -1.08%  '
Set-TextValue 'D42' '2.26'
Set-TextValue 'E42' '  -6.53%  '
Set-TextValue 'D43' '1.945.42'
Set-TextValue 'E43' '  -3.07%  '
Set-TextValue 'E44' '  -1.93%  '
Set-TextValue 'D45' '18.69'
Set-TextValue 'E45' '  -3.01%  '
Set-TextValue 'D46' '2.96'
Set-TextValue 'E46' '  -3.65%  '
Set-TextValue 'D47' '9.79'
Set-TextValue 'E47' '  -4.20%  '
Set-TextValue 'D48' '54.07'
Set-TextValue 'E48' '  +0.31%  '
Set-TextValue 'D49' '2.499.41'
Set-TextValue 'E49' '  -0.71%  '
Set-TextValue 'D50' '92.17'
Set-TextValue 'E50' '  -1.18%  '
Set-TextValue 'D51' '71.54'
Set-TextValue 'E51' '  -2.03%  '
